# Added course scheduling feature to project.
# The "Derslik" (classroom) table in this workbook had its "ID" column
# (the first column, A) removed - İsim / Tür / Kapasite now start at
# column A instead of B/C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column A first (mirrors the manual "right-click column header ->
# Delete" UI flow) so the resulting selection state matches a real delete.
$ws.Columns("A").Select() | Out-Null

# Remove the ID column entirely; this shifts İsim/Tür/Kapasite (and all
# their data) left by one column (B->A, C->B, D->C).
$ws.Columns("A").Delete()

# The worksheet's table ("Tablo1") still thinks it spans 4 columns -
# shrink it back down to the new A1:C11 extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C11"))

# Resizing the table doesn't refresh the ListColumns' cached names, so
# force that by re-writing the header cells with their own (now-shifted)
# values - this re-syncs "İsim"/"Tür"/"Kapasite" as the table's column
# names (instead of the stale "ID"/"İsim"/"Tür").
$ws.Range("A1").Value = $ws.Range("A1").Value()
$ws.Range("B1").Value = $ws.Range("B1").Value()
$ws.Range("C1").Value = $ws.Range("C1").Value()
